# Update "想去人数" (interested-count) values in column F on sheets
# "展览" (Exhibitions) and "全部类型" (All types), per the upstream
# data refresh captured in the commit "Update gh-pages to output
# generated at 456a3b4".

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# Sheet "展览" (sheet1): row -> new F value
$exhibitUpdates = @{
    2  = 7577
    5  = 4618
    8  = 610
    10 = 133
    11 = 432
    12 = 757
    13 = 26
    14 = 64
    15 = 247
    17 = 246
    18 = 130
    23 = 561
    25 = 681
    26 = 37
    27 = 35
    29 = 592
}

foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Range("F$row").Value = $exhibitUpdates[$row]
}

# Sheet "全部类型" (sheet4): row -> new F value
$allTypeUpdates = @{
    3  = 7577
    7  = 4618
    10 = 610
    13 = 133
    14 = 432
    18 = 757
    19 = 26
    20 = 64
    21 = 247
    26 = 246
    27 = 130
    32 = 561
    34 = 681
    35 = 37
    36 = 35
    38 = 592
}

foreach ($row in $allTypeUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allTypeUpdates[$row]
}
